$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text updates (Volume/Number + date range) ---
$ws.Range("A8").Value = "Volume 33   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/12/2026  Through  1/18/2026"

# --- Precinct data table updates (rows 15-30) ---

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 5
$ws.Range("I15").Value = 4
$ws.Range("N15").Value = 100

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 20
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -57.142857142857
$ws.Range("N16").Value = -93.181818181818

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 42.105263157894
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = -7.692307692307
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -68.421052631578

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = -70
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = -40
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -86.956521739130
$ws.Range("N18").Value = -97.619047619047

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 21
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 12
$ws.Range("K19").Value = 16.666666666666
$ws.Range("L19").Value = -6.666666666666
$ws.Range("M19").Value = -39.130434782608
$ws.Range("N19").Value = -17.647058823529

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 4
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -20
$ws.Range("M20").Value = -50
$ws.Range("N20").Value = -96

# Row 21
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 12.5
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = 3.076923076923
$ws.Range("I21").Value = 43
$ws.Range("J21").Value = 39
$ws.Range("K21").Value = 10.256410256410
$ws.Range("L21").Value = 4.878048780487
$ws.Range("M21").Value = -46.25
$ws.Range("N21").Value = -88.409703504043

# Row 22
$ws.Range("M22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M22").Value = -100

# Row 23
$ws.Range("L23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L23").Value = -100

# Row 24
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -44
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = -25.531914893617
$ws.Range("I24").Value = 41
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = -31.666666666666
$ws.Range("L24").Value = -6.818181818181
$ws.Range("M24").Value = -22.641509433962

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = -58.333333333333
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 13
$ws.Range("K25").Value = -46.153846153846
$ws.Range("L25").Value = -12.5

# Row 26
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 183.333333333333
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 67.857142857142
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 89.473684210526
$ws.Range("L26").Value = 56.521739130434
$ws.Range("M26").Value = 38.461538461538

# Row 27
$ws.Range("C27").Value = 1

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 2
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 66.666666666666
$ws.Range("L28").Value = 25

# Row 29
$ws.Range("N29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N29").Value = -100

# Row 30
$ws.Range("N30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N30").Value = -100

Write-Output "Edit applied successfully."
